$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '95.684.26'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.564.60'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +7.18%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.83'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '635.53'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.70%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.48'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +6.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.400'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("E10").Value = '  +8.67%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.561.70'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +7.16%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.45'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.33%  '

$ws.Range("E13").Value = '  +3.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.41'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +7.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.233.40'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +7.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '95.592.12'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.18%  '

$ws.Range("E17").Value = '  +4.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.557.65'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +7.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.17'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +20.72%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.99'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.84'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.500'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +10.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '517.06'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +4.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.41'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.96%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000194'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +7.00%  '

$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.68'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +8.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '92.88'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.66%  '

$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.26'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.78%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.07'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +17.62%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.146'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +5.61%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.49'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.93%  '

$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.03%  '

$ws.Range("B33").Value = 'Cronos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.184'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +6.05%  '

$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.41%  '

$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '30.10'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +6.27%  '

$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.565'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +6.65%  '

$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '579.15'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +9.34%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.86'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +6.02%  '

$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.48'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +8.27%  '

$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.02%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.151'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.05%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.924'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +6.91%  '

$ws.Range("B43").Value = 'ImmutableX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.75'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.29%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0431'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '23.81'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.02%  '

$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.61'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +4.00%  '

$ws.Range("B47").Value = 'MantraDAO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.54'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.44%  '

$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.17'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.19%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.74'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.61%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.16'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.37%  '

$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.11'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.45%  '

